$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 15573
$ws.Range("F6").Value = 419
$ws.Range("G6").Value = "已售罄"
$ws.Range("F8").Value = 703
$ws.Range("F9").Value = 15403
$ws.Range("F11").Value = 9004
$ws.Range("F15").Value = 89
$ws.Range("F19").Value = 18
$ws.Range("F29").Value = 83
$ws.Range("F34").Value = 40
$ws.Range("F35").Value = 250
$ws.Range("F36").Value = 317
$ws.Range("F37").Value = 451
$ws.Range("F39").Value = 5526

# --- Sheet "全部类型" (sheet4) ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F5").Value = 15573
$ws2.Range("F6").Value = 419
$ws2.Range("G6").Value = "已售罄"
$ws2.Range("F8").Value = 703
$ws2.Range("F9").Value = 15403
$ws2.Range("F11").Value = 9004
$ws2.Range("F15").Value = 89
$ws2.Range("F19").Value = 18
$ws2.Range("F29").Value = 83
$ws2.Range("F36").Value = 40
$ws2.Range("F37").Value = 250
$ws2.Range("F38").Value = 317
$ws2.Range("F39").Value = 451
$ws2.Range("F41").Value = 5526
